$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I18").Value = 2459
$ws.Range("K18").Value = 2459
$ws.Range("M18").Value = -2175
$ws.Range("H20").Value = 10300
$ws.Range("I20").Value = 10300
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 10300
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -10070
$ws.Range("N20").ClearContents()
$ws.Range("H28").Value = 1338.8636
$ws.Range("J28").Value = 3665
$ws.Range("L28").Value = 3665
$ws.Range("N28").Value = -4635
$ws.Range("H33").Value = 167.3
$ws.Range("I33").Value = 179
$ws.Range("J33").Value = 62
$ws.Range("K33").Value = 179
$ws.Range("L33").Value = 62
$ws.Range("M33").Value = 50
$ws.Range("N33").Value = -520
$ws.Range("H35").Value = 10300
$ws.Range("I35").Value = 10300
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 10300
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -9921
$ws.Range("N35").ClearContents()
$ws.Range("H55").Value = 820.6667
$ws.Range("J55").Value = 804.6667
$ws.Range("L55").Value = 804.6667
$ws.Range("N55").Value = -1232.6667
$ws.Range("H70").Value = 2939.5715
$ws.Range("J70").Value = 2935.4
$ws.Range("L70").Value = 8806.200000000001
$ws.Range("N70").Value = -9346.200000000001
$ws.Range("H73").Value = 2939.5715
$ws.Range("J73").Value = 2935.4
$ws.Range("L73").Value = 8806.200000000001
$ws.Range("N73").Value = -10678.2
$ws.Range("H106").Value = 1754.2727
$ws.Range("I106").Value = 1619.7
$ws.Range("K106").Value = 1619.7
$ws.Range("M106").Value = -988.7
$ws.Range("H132").Value = 5334.9116
$ws.Range("I132").Value = 5484.758
$ws.Range("K132").Value = 16454.274
$ws.Range("M132").Value = -13924.274
$ws.Range("H137").Value = 5168.92
$ws.Range("I137").Value = 1255.3684
$ws.Range("K137").Value = 3766.1052
$ws.Range("M137").Value = -1216.1052
$ws.Range("H138").Value = 311948.72
$ws.Range("I138").Value = 3727.2
$ws.Range("J138").Value = 449547.6
$ws.Range("K138").Value = 11181.6
$ws.Range("L138").Value = 1348642.8
$ws.Range("M138").Value = -6041.599999999999
$ws.Range("N138").Value = -1358922.8
$ws.Range("H141").Value = 7599.1875
$ws.Range("I141").Value = 4275.923
$ws.Range("K141").Value = 12827.769
$ws.Range("M141").Value = -7647.769

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4700.0713
$ws.Range("I32").Value = 4424.488
$ws.Range("K32").Value = 4424.488
$ws.Range("M32").Value = -4137.488
$ws.Range("H45").Value = 18677.785
$ws.Range("I45").Value = 24142.6
$ws.Range("J45").Value = 5015.75
$ws.Range("K45").Value = 24142.6
$ws.Range("L45").Value = 5015.75
$ws.Range("M45").Value = -23765.6
$ws.Range("N45").Value = -5769.75
$ws.Range("H61").Value = 3028.762
$ws.Range("I61").Value = 1639.2667
$ws.Range("K61").Value = 1639.2667
$ws.Range("M61").Value = -1427.2667
$ws.Range("H74").Value = 116948.06
$ws.Range("I74").Value = 143338.05
$ws.Range("K74").Value = 143338.05
$ws.Range("M74").Value = -142464.05
$ws.Range("H77").Value = 116948.06
$ws.Range("I77").Value = 143338.05
$ws.Range("K77").Value = 716690.25
$ws.Range("M77").Value = -712322.25
$ws.Range("H102").Value = 5731.7617
$ws.Range("I102").Value = 5637.933
$ws.Range("J102").Value = 5966.3335
$ws.Range("K102").Value = 5637.933
$ws.Range("L102").Value = 5966.3335
$ws.Range("M102").Value = -4015.933
$ws.Range("N102").Value = -9210.333500000001
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H136").Value = 3028.762
$ws.Range("I136").Value = 1639.2667
$ws.Range("K136").Value = 4917.800099999999
$ws.Range("M136").Value = -2367.800099999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3779
$ws.Range("J86").Value = 6500
$ws.Range("L86").Value = 6500
$ws.Range("N86").Value = -8746
$ws.Range("H89").Value = 3779
$ws.Range("J89").Value = 6500
$ws.Range("L89").Value = 32500
$ws.Range("N89").Value = -43732
$ws.Range("H105").Value = 18575000
$ws.Range("J105").Value = 27781876
$ws.Range("L105").Value = 27781876
$ws.Range("N105").Value = -27785370

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1399.186
$ws.Range("I16").Value = 1342.4333
$ws.Range("K16").Value = 1342.4333
$ws.Range("M16").Value = -1055.4333
$ws.Range("H31").Value = 6176.7896
$ws.Range("I31").Value = 4376.4614
$ws.Range("K31").Value = 4376.4614
$ws.Range("M31").Value = -4081.4614
$ws.Range("H34").Value = 6176.7896
$ws.Range("I34").Value = 4376.4614
$ws.Range("K34").Value = 4376.4614
$ws.Range("M34").Value = -4174.4614
$ws.Range("H58").Value = 2408.158
$ws.Range("I58").Value = 1656.4231
$ws.Range("K58").Value = 1656.4231
$ws.Range("M58").Value = -1453.4231
$ws.Range("H113").Value = 1399.186
$ws.Range("I113").Value = 1342.4333
$ws.Range("K113").Value = 1342.4333
$ws.Range("M113").Value = 827.5667000000001
$ws.Range("H122").Value = 2964.8125
$ws.Range("I122").Value = 2191.75
$ws.Range("K122").Value = 6575.25
$ws.Range("M122").Value = -4125.25
$ws.Range("H132").Value = 3470.3076
$ws.Range("I132").Value = 2868.375
$ws.Range("K132").Value = 8605.125
$ws.Range("M132").Value = -6075.125
$ws.Range("H136").Value = 2408.158
$ws.Range("I136").Value = 1656.4231
$ws.Range("K136").Value = 4969.2693
$ws.Range("M136").Value = -2419.2693

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1983.5238
$ws.Range("J5").Value = 2241.4443
$ws.Range("L5").Value = 6724.3329
$ws.Range("N5").Value = -6948.3329
$ws.Range("H28").Value = 179
$ws.Range("I28").Value = 179
$ws.Range("K28").Value = 537
$ws.Range("M28").Value = -305
$ws.Range("H49").Value = 2500
$ws.Range("J49").Value = 2500
$ws.Range("L49").Value = 7500
$ws.Range("N49").Value = -7812
$ws.Range("H113").Value = 5473
$ws.Range("I113").Value = 417.625
$ws.Range("J113").Value = 7852
$ws.Range("K113").Value = 1252.875
$ws.Range("L113").Value = 23556
$ws.Range("M113").Value = 917.125
$ws.Range("N113").Value = -27896
$ws.Range("H135").Value = 1983.5238
$ws.Range("J135").Value = 2241.4443
$ws.Range("L135").Value = 20172.9987
$ws.Range("N135").Value = -25242.9987

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6915.4243
$ws.Range("I122").Value = 6608.52
$ws.Range("J122").Value = 7874.5
$ws.Range("K122").Value = 19825.56
$ws.Range("L122").Value = 23623.5
$ws.Range("M122").Value = -17375.56
$ws.Range("N122").Value = -28523.5
$ws.Range("H126").Value = 7613.3335
$ws.Range("I126").Value = 1904.4
$ws.Range("K126").Value = 5713.200000000001
$ws.Range("M126").Value = -3243.200000000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4475.125
$ws.Range("I61").Value = 4400.143
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 4400.143
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -4198.143
$ws.Range("N61").Value = -5404
$ws.Range("H68").Value = 2674.9
$ws.Range("I68").Value = 2358.8333
$ws.Range("J68").Value = 3149
$ws.Range("K68").Value = 2358.8333
$ws.Range("L68").Value = 3149
$ws.Range("M68").Value = -1609.8333
$ws.Range("N68").Value = -4647
$ws.Range("H71").Value = 2674.9
$ws.Range("I71").Value = 2358.8333
$ws.Range("J71").Value = 3149
$ws.Range("K71").Value = 11794.1665
$ws.Range("L71").Value = 15745
$ws.Range("M71").Value = -8050.166499999999
$ws.Range("N71").Value = -23233
$ws.Range("H113").Value = 4475.125
$ws.Range("I113").Value = 4400.143
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 4400.143
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -2230.143
$ws.Range("N113").Value = -9340
$ws.Range("H136").Value = 3413.2
$ws.Range("I136").Value = 3207.2307
$ws.Range("J136").Value = 3795.7144
$ws.Range("K136").Value = 9621.6921
$ws.Range("L136").Value = 11387.1432
$ws.Range("M136").Value = -7071.6921
$ws.Range("N136").Value = -16487.1432

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4936.8276
$ws.Range("I62").Value = 4106.577
$ws.Range("K62").Value = 4106.577
$ws.Range("M62").Value = -3482.577
$ws.Range("H65").Value = 4936.8276
$ws.Range("I65").Value = 4106.577
$ws.Range("K65").Value = 20532.885
$ws.Range("M65").Value = -17412.885
$ws.Range("H113").Value = 549.5217
$ws.Range("J113").Value = 610.75
$ws.Range("L113").Value = 1832.25
$ws.Range("N113").Value = -6172.25
$ws.Range("H122").Value = 15629946
$ws.Range("I122").Value = 5947.8184
$ws.Range("J122").Value = 50002740
$ws.Range("K122").Value = 17843.4552
$ws.Range("L122").Value = 150008220
$ws.Range("M122").Value = -15393.4552
$ws.Range("N122").Value = -150013120
$ws.Range("H132").Value = 3391.9355
$ws.Range("J132").Value = 2843.5557
$ws.Range("L132").Value = 8530.667099999999
$ws.Range("N132").Value = -13590.6671
$ws.Range("H136").Value = 41670108
$ws.Range("I136").Value = 52632776
$ws.Range("K136").Value = 157898328
$ws.Range("M136").Value = -157895778
